$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = "8.300,01 TL - 199,41 TL"
$ws.Range("H6").Value = ""

$ws.Range("D12").Value = "WU: 0,75 USD–12 USD; Diğer: 700 TL–4.000 TL"

$ws.Range("D13").Value = "Hesaba: Asgari 1 TL | Azami 909,5 TL"
$ws.Range("F13").Value = ""
$ws.Range("H13").Value = ""

$ws.Range("D14").Value = "3.500 TL - 13.500 TL"
$ws.Range("F14").Value = ""
$ws.Range("H14").Value = ""
